$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: the AlternateContent x15ac:absPath metadata (internal SAS output-folder path
# bookkeeping) is not exposed via the Excel object model, so it is left untouched.

# Update report title (Year-to-Date period, October -> November 2016)
$ws.Range("A2").Value = "Industrial Sector by Census Division and State, Year-to-Date through November 2016"

# Update the data table values (revised Relative Standard Error figures)
$ws.Cells.Item(4, 2).Value = 81
$ws.Cells.Item(4, 3).Value = 141
$ws.Cells.Item(4, 5).Value = 36
$ws.Cells.Item(4, 8).Value = 47
$ws.Cells.Item(5, 3).Value = 436
$ws.Cells.Item(5, 5).Value = 55
$ws.Cells.Item(6, 3).Value = 128
$ws.Cells.Item(6, 5).Value = 51
$ws.Cells.Item(6, 8).Value = 47
$ws.Cells.Item(7, 2).Value = 195
$ws.Cells.Item(7, 3).Value = 4303
$ws.Cells.Item(7, 8).Value = 788
$ws.Cells.Item(8, 3).Value = 418
$ws.Cells.Item(8, 5).Value = 209
$ws.Cells.Item(9, 2).Value = 22
$ws.Cells.Item(9, 4).Value = 69
$ws.Cells.Item(9, 5).Value = 24
$ws.Cells.Item(9, 6).Value = 29
$ws.Cells.Item(9, 8).Value = 182
$ws.Cells.Item(10, 3).Value = 529
$ws.Cells.Item(10, 4).Value = 128
$ws.Cells.Item(10, 5).Value = 64
$ws.Cells.Item(10, 6).Value = 75
$ws.Cells.Item(11, 3).Value = 46
$ws.Cells.Item(11, 5).Value = 38
$ws.Cells.Item(11, 8).Value = 182
$ws.Cells.Item(12, 2).Value = 47
$ws.Cells.Item(12, 3).Value = 25
$ws.Cells.Item(12, 4).Value = 81
$ws.Cells.Item(12, 5).Value = 31
$ws.Cells.Item(12, 6).Value = 28
$ws.Cells.Item(13, 2).Value = 9
$ws.Cells.Item(13, 3).Value = 23
$ws.Cells.Item(13, 4).Value = 78
$ws.Cells.Item(13, 5).Value = 18
$ws.Cells.Item(13, 6).Value = 21
$ws.Cells.Item(13, 8).Value = 82
$ws.Cells.Item(14, 2).Value = 8
$ws.Cells.Item(14, 5).Value = 50
$ws.Cells.Item(14, 6).Value = 99
$ws.Cells.Item(15, 2).Value = 614
$ws.Cells.Item(15, 3).Value = 8
$ws.Cells.Item(15, 5).Value = 32
$ws.Cells.Item(15, 6).Value = 16
$ws.Cells.Item(16, 2).Value = 99
$ws.Cells.Item(16, 3).Value = 19
$ws.Cells.Item(16, 4).Value = 89
$ws.Cells.Item(16, 5).Value = 33
$ws.Cells.Item(16, 8).Value = 208
$ws.Cells.Item(17, 2).Value = 125
$ws.Cells.Item(17, 3).Value = 213
$ws.Cells.Item(17, 4).Value = 339
$ws.Cells.Item(17, 5).Value = 43
$ws.Cells.Item(17, 6).Value = 168
$ws.Cells.Item(18, 2).Value = 18
$ws.Cells.Item(18, 3).Value = 374
$ws.Cells.Item(18, 5).Value = 37
$ws.Cells.Item(18, 8).Value = 89
$ws.Cells.Item(19, 2).Value = 15
$ws.Cells.Item(19, 3).Value = 224
$ws.Cells.Item(19, 4).Value = 168
$ws.Cells.Item(19, 5).Value = 29
$ws.Cells.Item(19, 6).Value = 100
$ws.Cells.Item(19, 8).Value = 97
$ws.Cells.Item(20, 2).Value = 15
$ws.Cells.Item(20, 3).Value = 429
$ws.Cells.Item(20, 4).Value = 168
$ws.Cells.Item(20, 5).Value = 30
$ws.Cells.Item(21, 5).Value = 112
$ws.Cells.Item(22, 2).Value = 37
$ws.Cells.Item(22, 3).Value = 409
$ws.Cells.Item(22, 5).Value = 87
$ws.Cells.Item(22, 8).Value = 97
$ws.Cells.Item(23, 2).Value = 175
$ws.Cells.Item(23, 5).Value = 306
$ws.Cells.Item(24, 2).Value = 44
$ws.Cells.Item(24, 5).Value = 457
$ws.Cells.Item(25, 2).Value = 90
$ws.Cells.Item(25, 3).Value = 290
$ws.Cells.Item(25, 5).Value = 175
$ws.Cells.Item(25, 6).Value = 100
$ws.Cells.Item(26, 2).Value = 21
$ws.Cells.Item(26, 3).Value = 74
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 8).Value = 46
$ws.Cells.Item(28, 2).Value = 85
$ws.Cells.Item(28, 3).Value = 201
$ws.Cells.Item(29, 2).Value = 42
$ws.Cells.Item(29, 3).Value = 78
$ws.Cells.Item(29, 5).Value = 25
$ws.Cells.Item(29, 8).Value = 288
$ws.Cells.Item(30, 3).Value = 277
$ws.Cells.Item(30, 5).Value = 87
$ws.Cells.Item(31, 2).Value = 82
$ws.Cells.Item(31, 3).Value = 666
$ws.Cells.Item(31, 5).Value = 70
$ws.Cells.Item(31, 8).Value = 1086
$ws.Cells.Item(32, 2).Value = 9
$ws.Cells.Item(32, 3).Value = 8
$ws.Cells.Item(32, 5).Value = 75
$ws.Cells.Item(33, 2).Value = 25
$ws.Cells.Item(33, 3).Value = 570
$ws.Cells.Item(33, 5).Value = 29
$ws.Cells.Item(33, 8).Value = 405
$ws.Cells.Item(34, 8).Value = 30
$ws.Cells.Item(35, 2).Value = 6
$ws.Cells.Item(35, 3).Value = 88
$ws.Cells.Item(35, 5).Value = 15
$ws.Cells.Item(35, 6).Value = 72
$ws.Cells.Item(36, 2).Value = 65
$ws.Cells.Item(36, 3).Value = 104
$ws.Cells.Item(36, 5).Value = 22
$ws.Cells.Item(36, 6).Value = 119
$ws.Cells.Item(37, 5).Value = 71
$ws.Cells.Item(38, 5).Value = 38
$ws.Cells.Item(39, 3).Value = 148
$ws.Cells.Item(40, 2).Value = 42
$ws.Cells.Item(40, 3).Value = 67
$ws.Cells.Item(40, 4).Value = 55
$ws.Cells.Item(40, 6).Value = 7
$ws.Cells.Item(41, 3).Value = 49
$ws.Cells.Item(41, 5).Value = 24
$ws.Cells.Item(42, 4).Value = 78
$ws.Cells.Item(42, 6).Value = 7
$ws.Cells.Item(43, 2).Value = 50
$ws.Cells.Item(43, 3).Value = 88
$ws.Cells.Item(43, 5).Value = 84
$ws.Cells.Item(44, 3).Value = 408
$ws.Cells.Item(44, 4).Value = 59
$ws.Cells.Item(44, 6).Value = 12
$ws.Cells.Item(45, 2).Value = 42
$ws.Cells.Item(45, 3).Value = 743
$ws.Cells.Item(45, 6).Value = 9
$ws.Cells.Item(46, 2).Value = 370
$ws.Cells.Item(46, 3).Value = 678
$ws.Cells.Item(46, 5).Value = 80
$ws.Cells.Item(47, 2).Value = 105
$ws.Cells.Item(47, 5).Value = 50
$ws.Cells.Item(48, 2).Value = 308
$ws.Cells.Item(49, 5).Value = 21
$ws.Cells.Item(50, 3).Value = 3020
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 3).Value = 1341
$ws.Cells.Item(51, 5).Value = 16
$ws.Cells.Item(51, 6).Value = 457
$ws.Cells.Item(52, 2).Value = 46
$ws.Cells.Item(52, 3).Value = 372
$ws.Cells.Item(52, 5).Value = 16
$ws.Cells.Item(52, 6).Value = 7
$ws.Cells.Item(53, 3).Value = 95
$ws.Cells.Item(53, 5).Value = 3
$ws.Cells.Item(53, 6).Value = 8
$ws.Cells.Item(54, 3).Value = 405
$ws.Cells.Item(54, 5).Value = 3
$ws.Cells.Item(54, 6).Value = 8
$ws.Cells.Item(55, 5).Value = 72
$ws.Cells.Item(56, 3).Value = 82
$ws.Cells.Item(57, 2).Value = 201
$ws.Cells.Item(57, 3).Value = 53
$ws.Cells.Item(57, 5).Value = 106
$ws.Cells.Item(57, 6).Value = 130
$ws.Cells.Item(57, 8).Value = 149
$ws.Cells.Item(58, 3).Value = 19
$ws.Cells.Item(58, 5).Value = 106
$ws.Cells.Item(59, 2).Value = 201
$ws.Cells.Item(59, 3).Value = 62
$ws.Cells.Item(59, 6).Value = 130
$ws.Cells.Item(59, 8).Value = 149
$ws.Cells.Item(60, 2).Value = 7
$ws.Cells.Item(60, 3).Value = 33
$ws.Cells.Item(60, 4).Value = 34
$ws.Cells.Item(60, 6).Value = 7
$ws.Cells.Item(60, 8).Value = 32
